$d = $word.ActiveDocument

# 1. Update hyperlink field instrText (filename change)
$d.Content.Find.Execute(
    "certificate-of-completion-connor-readnour.pdf",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "javascript-fundamentals-connor-readnour.pdf", 2
)

Write-Output "done"
